# Updated cryptos list on Mon Jun  5 13:34:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the Price column (D) while forcing it to stay
# stored as text (many of the price strings look numeric, e.g. "1.001", and
# Excel would otherwise silently convert them to numbers on assignment).
function Set-PriceText($row, $text) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows whose Price (D) and Volume(1h) (E) values both changed.
$data = @{
    2 = @("26.738.55", "  -1.87%  ")
    3 = @("1.868.18", "  -2.04%  ")
    5 = @("300.17", "  -2.34%  ")
    6 = @("1.000", "  -0.03%  ")
    7 = @("0.5321", "  +1.10%  ")
    8 = @("0.3728", "  -2.24%  ")
    9 = @("0.07158", "  -1.71%  ")
    10 = @("21.46", "  -1.99%  ")
    11 = @("0.8860", "  -1.81%  ")
    12 = @("0.08172", "  +0.10%  ")
    13 = @("1.868.14", "  +27.39%  ")
    14 = @("92.24", "  -4.13%  ")
    15 = @("5.288", "  -1.38%  ")
    16 = @("1.001", "  -0.03%  ")
    17 = @("14.82", "  +0.38%  ")
    18 = @("0.000008480", "  -2.03%  ")
    20 = @("26.777.95", "  -1.86%  ")
    21 = @("4.969", "  -2.92%  ")
    22 = @("10.61", "  -2.03%  ")
    23 = @("6.357", "  -2.48%  ")
    24 = @("2.287", "  -1.01%  ")
    25 = @("145.58", "  -2.92%  ")
    28 = @("113.61", "  -2.71%  ")
    29 = @("4.686", "  -3.34%  ")
    30 = @("4.623", "  -4.71%  ")
    31 = @("0.09109", "  -1.47%  ")
    32 = @("0.8033", "  -3.31%  ")
    33 = @("0.05008", "  -1.15%  ")
    34 = @("1.169", "  -4.79%  ")
    35 = @("2.940", "  -1.61%  ")
    36 = @("0.6099", "  +4.96%  ")
    38 = @("3.173", "  -5.24%  ")
    39 = @("0.01940", "  -3.13%  ")
    40 = @("1.061", "  -1.66%  ")
    43 = @("8.715", "  -4.86%  ")
    44 = @("114.74", "  -1.73%  ")
    45 = @("0.1491", "  -2.10%  ")
    46 = @("1.000", "  -0.02%  ")
    47 = @("1.636", "  -0.58%  ")
    48 = @("9.890", "  -3.15%  ")
    49 = @("37.28", "  -4.39%  ")
    50 = @("0.06058", "  -1.11%  ")
    51 = @("62.03", "  -3.92%  ")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    Set-PriceText $row $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Rows where only Volume(1h) (E) changed.
$eOnly = @{
    19 = "  -0.02%  "
    26 = "  -0.56%  "
    37 = "  -1.97%  "
}

foreach ($row in $eOnly.Keys) {
    $ws.Cells.Item($row, 5).Value = $eOnly[$row]
}

# Rows 41 and 42 swapped contents (Decentraland moved above FraxShare)
# with refreshed Price / Volume(1h) values.
$ws.Cells.Item(41, 2).Value = "Decentraland"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-PriceText 41 "0.5213"
$ws.Cells.Item(41, 5).Value = "  +5.74%  "

$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-PriceText 42 "6.476"
$ws.Cells.Item(42, 5).Value = "  -1.87%  "
